$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 231.125
$ws.Range("I9").Value = 245
$ws.Range("J9").Value = 189.5
$ws.Range("K9").Value = 245
$ws.Range("L9").Value = 189.5
$ws.Range("M9").Value = -76
$ws.Range("N9").Value = -527.5
$ws.Range("H12").Value = 280.5
$ws.Range("I12").Value = 306
$ws.Range("J12").Value = 102
$ws.Range("K12").Value = 306
$ws.Range("L12").Value = 102
$ws.Range("M12").Value = -136
$ws.Range("N12").Value = -442
$ws.Range("H42").Value = 434.53845
$ws.Range("I42").Value = 143.16667
$ws.Range("J42").Value = 684.2857
$ws.Range("K42").Value = 429.50001
$ws.Range("L42").Value = 2052.8571
$ws.Range("M42").Value = -199.50001
$ws.Range("N42").Value = -2512.8571
$ws.Range("H53").Value = 154.63158
$ws.Range("I53").Value = 170
$ws.Range("J53").Value = 147.53847
$ws.Range("K53").Value = 170
$ws.Range("L53").Value = 147.53847
$ws.Range("M53").Value = 467
$ws.Range("N53").Value = -1421.53847
$ws.Range("H70").Value = 1350.1
$ws.Range("I70").Value = 1220.2
$ws.Range("J70").Value = 1480
$ws.Range("K70").Value = 3660.6
$ws.Range("L70").Value = 4440
$ws.Range("M70").Value = -3390.6
$ws.Range("N70").Value = -4980
$ws.Range("H73").Value = 1350.1
$ws.Range("I73").Value = 1220.2
$ws.Range("J73").Value = 1480
$ws.Range("K73").Value = 3660.6
$ws.Range("L73").Value = 4440
$ws.Range("M73").Value = -2724.6
$ws.Range("N73").Value = -6312
$ws.Range("H100").Value = 1555.7693
$ws.Range("I100").Value = 996.3333
$ws.Range("J100").Value = 2035.2858
$ws.Range("K100").Value = 996.3333
$ws.Range("L100").Value = 2035.2858
$ws.Range("M100").Value = -455.3333
$ws.Range("N100").Value = -3117.2858
$ws.Range("H116").Value = 2971.7334
$ws.Range("J116").Value = 3668.6667
$ws.Range("L116").Value = 3668.6667
$ws.Range("N116").Value = -10552.6667
$ws.Range("H132").Value = 7998.727
$ws.Range("I132").Value = 8998.375
$ws.Range("K132").Value = 26995.125
$ws.Range("M132").Value = -24465.125
$ws.Range("H138").Value = 2994.8
$ws.Range("I138").Value = 2613.6538
$ws.Range("J138").Value = 3248.8975
$ws.Range("K138").Value = 7840.9614
$ws.Range("L138").Value = 9746.692500000001
$ws.Range("M138").Value = -2700.9614
$ws.Range("N138").Value = -20026.6925
$ws.Range("H141").Value = 4437.1665
$ws.Range("I141").Value = 1737.9
$ws.Range("J141").Value = 7811.25
$ws.Range("K141").Value = 5213.700000000001
$ws.Range("L141").Value = 23433.75
$ws.Range("M141").Value = -33.70000000000073
$ws.Range("N141").Value = -33793.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 18103
$ws.Range("I31").Value = 1575.7142
$ws.Range("K31").Value = 1575.7142
$ws.Range("M31").Value = -1281.7142
$ws.Range("H61").Value = 9807172
$ws.Range("I61").Value = 17546576
$ws.Range("K61").Value = 17546576
$ws.Range("M61").Value = -17546364
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null
$ws.Range("H123").Value = 39571.6
$ws.Range("J123").Value = 39571.6
$ws.Range("L123").Value = 39571.6
$ws.Range("N123").Value = -49371.6
$ws.Range("H136").Value = 9807172
$ws.Range("I136").Value = 17546576
$ws.Range("K136").Value = 52639728
$ws.Range("M136").Value = -52637178
$ws.Range("H139").Value = 71211
$ws.Range("J139").Value = 69123.336
$ws.Range("L139").Value = 69123.336
$ws.Range("N139").Value = -79403.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1727.2
$ws.Range("I99").Value = 1426.5385
$ws.Range("J99").Value = 2285.5715
$ws.Range("K99").Value = 1426.5385
$ws.Range("L99").Value = 2285.5715
$ws.Range("M99").Value = 71.46149999999989
$ws.Range("N99").Value = -5281.5715
$ws.Range("H102").Value = 39907.855
$ws.Range("I102").Value = 18071
$ws.Range("J102").Value = 94500
$ws.Range("K102").Value = 18071
$ws.Range("L102").Value = 94500
$ws.Range("M102").Value = -14826
$ws.Range("N102").Value = -100990
$ws.Range("H134").Value = 3078.4348
$ws.Range("I134").Value = 2752.9412
$ws.Range("J134").Value = 4000.6667
$ws.Range("K134").Value = 8258.8236
$ws.Range("L134").Value = 12002.0001
$ws.Range("M134").Value = -5723.8236
$ws.Range("N134").Value = -17072.0001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6650.2617
$ws.Range("I31").Value = 1962.1428
$ws.Range("J31").Value = 7587.8857
$ws.Range("K31").Value = 1962.1428
$ws.Range("L31").Value = 7587.8857
$ws.Range("M31").Value = -1667.1428
$ws.Range("N31").Value = -8177.8857
$ws.Range("H34").Value = 6650.2617
$ws.Range("I34").Value = 1962.1428
$ws.Range("J34").Value = 7587.8857
$ws.Range("K34").Value = 1962.1428
$ws.Range("L34").Value = 7587.8857
$ws.Range("M34").Value = -1760.1428
$ws.Range("N34").Value = -7991.8857
$ws.Range("H58").Value = 1832.5082
$ws.Range("I58").Value = 1489.1842
$ws.Range("J58").Value = 2399.739
$ws.Range("K58").Value = 1489.1842
$ws.Range("L58").Value = 2399.739
$ws.Range("M58").Value = -1286.1842
$ws.Range("N58").Value = -2805.739
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H136").Value = 1832.5082
$ws.Range("I136").Value = 1489.1842
$ws.Range("J136").Value = 2399.739
$ws.Range("K136").Value = 4467.5526
$ws.Range("L136").Value = 7199.217000000001
$ws.Range("M136").Value = -1917.5526
$ws.Range("N136").Value = -12299.217

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5689.5713
$ws.Range("J131").Value = 6839.478
$ws.Range("L131").Value = 20518.434
$ws.Range("N131").Value = -30598.434

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2305.4285
$ws.Range("I126").Value = 2439.6667
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 7319.000100000001
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -4849.000100000001
$ws.Range("N126").Value = -9440

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 37777.668
$ws.Range("J115").Value = 37777.668
$ws.Range("L115").Value = 37777.668
$ws.Range("N115").Value = -40127.668
$ws.Range("H132").Value = 3490.15
$ws.Range("I132").Value = 2485.1538
$ws.Range("J132").Value = 5356.5713
$ws.Range("K132").Value = 7455.4614
$ws.Range("L132").Value = 16069.7139
$ws.Range("M132").Value = -4925.4614
$ws.Range("N132").Value = -21129.7139
$ws.Range("H133").Value = 44108.668
$ws.Range("J133").Value = 44108.668
$ws.Range("L133").Value = 44108.668
$ws.Range("N133").Value = -49168.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -15298
$ws.Range("H132").Value = 6785790
$ws.Range("I132").Value = 2846.1924
$ws.Range("J132").Value = 17159704
$ws.Range("K132").Value = 8538.5772
$ws.Range("L132").Value = 51479112
$ws.Range("M132").Value = -6008.5772
$ws.Range("N132").Value = -51484172
